$d = $word.ActiveDocument

# 1) "It would be a more complex algorithm" -> "In practice, this would be a
#    more complex algorithm" (in the second list item, right after
#    "...identify comorbidities that occur with age and behavior.")
$d.Content.Find.Execute(
    "identify comorbidities that occur with age and behavior. It would be",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "identify comorbidities that occur with age and behavior. In practice, this would be",
    2)

# 2) Add a new (un-numbered) list paragraph right after the second bullet
#    ("...Clinical Decision Support and Population Health.") and before the
#    third bullet ("Severable wearable devices ...").
$thirdBullet = $d.Paragraphs.Item(3)
$insertionPoint = $thirdBullet.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(3)
$newPara.Range.ListFormat.RemoveNumbers()

$newText = "For the project, we would want to pull from population health studies including population genetics and compile the data. If it were implemented in reality, we would want to aggregate the same data and insurance database data from many sources."

# Insert the paragraph text plus a one-character sentinel so the later
# collapsed-bookmark insertion point isn't the very last character before
# the paragraph mark (avoids a degenerate boundary position), then strip
# the sentinel back out again.
$textRange = $newPara.Range
$textRange.Collapse(1)
$textRange.InsertAfter($newText + "#")

$sentinelPos = $d.Paragraphs.Item(3).Range.End - 2
$bookmarkRange = $d.Range($sentinelPos, $sentinelPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$sentinelRange = $d.Range($sentinelPos, $sentinelPos + 1)
$sentinelRange.Delete()
